# Regenerate save_data to use K (strikeouts) instead of Strike# values,
# writing the recalculated s_vals into column G (header "K") for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2-19 (column G), replacing the old Strike# based values.
$kValues = @{
    2  = 3
    3  = 1
    4  = 2
    5  = 1
    6  = 4
    7  = 4
    8  = 3
    9  = 2
    10 = 0
    11 = 0
    12 = 1
    13 = 0
    14 = 0
    15 = 4
    16 = 2
    17 = 2
    18 = 0
    19 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
